$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(472).Insert()

$ws.Range("A472").Value = 3
$ws.Range("B472").Value = "Femacal de La Calera"
$ws.Range("C472").Value = "Coquimbo"
$ws.Range("D472").Value = 45135
$ws.Range("E472").Value = 5
$ws.Range("F472").Value = 100112009
$ws.Range("G472").Value = "Acelga"
$ws.Range("H472").Value = "Sin especificar"
$ws.Range("I472").Value = "Primera"
$ws.Range("J472").Value = 120
$ws.Range("K472").Value = 3000
$ws.Range("L472").Value = 3000
$ws.Range("M472").Value = 3000
$ws.Range("N472").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O472").Value = "Provincia de Quillota"
$ws.Range("P472").Value = 500
$ws.Range("Q472").Value = 6
$ws.Range("R472").Value = "Hortaliza"
